# Update countries & provincias Spain
# Applies the daily data refresh: new case counts for several countries,
# including Turquia and Burkina Faso overtaking their neighbours in the
# "Casos totales" ranking, plus the refreshed "datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 19:20"

# --- Estados Unidos (row 4) : new numbers, same rank -------------------
$ws.Range("B4").Value = 176518
$ws.Range("C4").Value = 12730
$ws.Range("D4").Value = 6241
$ws.Range("E4").Value = 166846
$ws.Range("F4").Value = 3893
$ws.Range("G4").Value = 290
$ws.Range("H4").Value = 3431

# --- Turquia overtakes Belgica and Paises Bajos (rows 13-15) -----------
# Row 13: Turquia (new data, now ranked ahead of Belgica / Paises Bajos)
$ws.Range("A13").Value = "Turquia"
$ws.Range("B13").Value = 13531
$ws.Range("C13").Value = 2704
$ws.Range("D13").Value = 243
$ws.Range("E13").Value = 13074
$ws.Range("F13").Value = 568
$ws.Range("G13").Value = 46
$ws.Range("H13").Value = 214

# Row 14: Belgica (values unchanged, just pushed down one rank)
$ws.Range("A14").Value = "Belgica"
$ws.Range("B14").Value = 12775
$ws.Range("C14").Value = 876
$ws.Range("D14").Value = 1696
$ws.Range("E14").Value = 10374
$ws.Range("F14").Value = 1021
$ws.Range("G14").Value = 192
$ws.Range("H14").Value = 705

# Row 15: Paises Bajos (values unchanged, just pushed down one rank)
$ws.Range("A15").Value = "Paises Bajos"
$ws.Range("B15").Value = 12595
$ws.Range("C15").Value = 845
$ws.Range("D15").Value = 250
$ws.Range("E15").Value = 11306
$ws.Range("F15").Value = 1053
$ws.Range("G15").Value = 175
$ws.Range("H15").Value = 1039

# --- Austria (row 16) : new numbers, same rank --------------------------
$ws.Range("B16").Value = 10109
$ws.Range("C16").Value = 491
$ws.Range("D16").Value = 1095
$ws.Range("E16").Value = 8886
$ws.Range("F16").Value = 198
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = 128

# --- Chequia (row 25) : new numbers, same rank --------------------------
$ws.Range("B25").Value = 3257
$ws.Range("C25").Value = 256
$ws.Range("D25").Value = 25
$ws.Range("E25").Value = 3207
$ws.Range("F25").Value = 64
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 25

# --- Marruecos (row 65) : new numbers, same rank -------------------------
$ws.Range("B65").Value = 602
$ws.Range("C65").Value = 46
$ws.Range("D65").Value = 24
$ws.Range("E65").Value = 542
$ws.Range("F65").Value = 1
$ws.Range("G65").Value = 3
$ws.Range("H65").Value = 36

# --- Lituania (row 68) : new numbers, same rank ---------------------------
$ws.Range("B68").Value = 537
$ws.Range("C68").Value = 46
$ws.Range("D68").Value = 7
$ws.Range("E68").Value = 522
$ws.Range("F68").Value = 27
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 8

# --- Jordania (row 86) : new numbers, same rank ---------------------------
$ws.Range("B86").Value = 274
$ws.Range("C86").Value = 6
$ws.Range("D86").Value = 30
$ws.Range("E86").Value = 239
$ws.Range("F86").Value = 5
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 5

# --- Burkina Faso overtakes Reunion (rows 88-89) --------------------------
# Row 88: Burkina Faso (new data, now ranked ahead of Reunion)
$ws.Range("A88").Value = "Burkina Faso"
$ws.Range("B88").Value = 261
$ws.Range("C88").Value = 15
$ws.Range("D88").Value = 32
$ws.Range("E88").Value = 215
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 2
$ws.Range("H88").Value = 14

# Row 89: Reunion (values unchanged, just pushed down one rank)
$ws.Range("A89").Value = "Reunion"
$ws.Range("B89").Value = 247
$ws.Range("C89").Value = 23
$ws.Range("D89").Value = 1
$ws.Range("E89").Value = 246
$ws.Range("F89").Value = 4
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 0

# --- Albania (row 90) : new numbers, same rank ------------------------------
$ws.Range("B90").Value = 243
$ws.Range("C90").Value = 20
$ws.Range("D90").Value = 52
$ws.Range("E90").Value = 176
$ws.Range("F90").Value = 8
$ws.Range("G90").Value = 4
$ws.Range("H90").Value = 15
